# SMART-3Seq-HT_barcodes.xlsx — "Add files via upload"
#
# Adds 20 new P5-PCR index rows (rows 6-25) to the "P5-PCR" sheet, with
# accompanying Name / Barcode / Reverse-complement-sequence data, widens
# column B to fit the longer names, and restores the cursor/selection on
# each of the three sheets to where the author last left it.

$wb = $excel.ActiveWorkbook

$wsOligo = $wb.Worksheets.Item("oligo-dT")
$wsTso   = $wb.Worksheets.Item("TSO")
$wsP5    = $wb.Worksheets.Item("P5-PCR")

# ---------------------------------------------------------------------
# New P5-PCR index data (rows 6-25 => index numbers 5-24)
# ---------------------------------------------------------------------
$names = @(
    "oNM1874_PCR_P5-index#5_AGGTGCGT",
    "oNM1875_PCR_P5-index#6_GAACATAC",
    "oNM1876_PCR_P5-index#7_ACATAGCG",
    "oNM1877_PCR_P5-index#8_GTGCGATA",
    "oNM1878_PCR_P5-index#9_CCAACAGA",
    "oNM1879_PCR_P5-index#10_TTGGTGAG",
    "oNM1880_PCR_P5-index#11_CGCGGTTC",
    "oNM1881_PCR_P5-index#12_TATAACCT",
    "oNM1882_PCR_P5-index#13_AAGGATGA",
    "oNM1883_PCR_P5-index#14_GGAAGCAG",
    "oNM1884_PCR_P5-index#15_TGACGAAT",
    "oNM1885_PCR_P5-index#16_CAGTAGGC",
    "oNM1886_PCR_P5-index#17_ATATTCAC",
    "oNM1887_PCR_P5-index#18_GCGCCTGT",
    "oNM1888_PCR_P5-index#19_ACTCTATG",
    "oNM1889_PCR_P5-index#20_GTCTCGCA",
    "oNM1890_PCR_P5-index#21_AAGACGTC",
    "oNM1891_PCR_P5-index#22_GGAGTACT",
    "oNM1892_PCR_P5-index#23_ACCGGCCA",
    "oNM1893_PCR_P5-index#24_GTTAATTG"
)

$barcodes = @(
    "AGGTGCGT",
    "GAACATAC",
    "ACATAGCG",
    "GTGCGATA",
    "CCAACAGA",
    "TTGGTGAG",
    "CGCGGTTC",
    "TATAACCT",
    "AAGGATGA",
    "GGAAGCAG",
    "TGACGAAT",
    "CAGTAGGC",
    "ATATTCAC",
    "GCGCCTGT",
    "ACTCTATG",
    "GTCTCGCA",
    "AAGACGTC",
    "GGAGTACT",
    "ACCGGCCA",
    "GTTAATTG"
)

$sequences = @(
    "AATGATACGGCGACCACCGAGATCTACACAGGTGCGTACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACGAACATACACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACACATAGCGACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACGTGCGATAACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACCCAACAGAACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACTTGGTGAGACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACCGCGGTTCACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACTATAACCTACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACAAGGATGAACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACGGAAGCAGACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACTGACGAATACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACCAGTAGGCACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACATATTCACACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACGCGCCTGTACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACACTCTATGACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACGTCTCGCAACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACAAGACGTCACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACGGAGTACTACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACACCGGCCAACACTCTTTCCCTACACGACGCTCTTCCGAT*C",
    "AATGATACGGCGACCACCGAGATCTACACGTTAATTGACACTCTTTCCCTACACGACGCTCTTCCGAT*C"
)

$firstRow = 6
$lastRow  = 25
$rowCount = $lastRow - $firstRow + 1

# Column A: sequential index numbers 5..24
for ($i = 0; $i -lt $rowCount; $i++) {
    $wsP5.Cells.Item($firstRow + $i, 1).Value = $i + 5
}

# Column B: names, in row order
for ($i = 0; $i -lt $rowCount; $i++) {
    $wsP5.Cells.Item($firstRow + $i, 2).Value = $names[$i]
}

# Column D: reverse-complement sequences, in row order
for ($i = 0; $i -lt $rowCount; $i++) {
    $wsP5.Cells.Item($firstRow + $i, 4).Value = $sequences[$i]
}

# Column C: barcodes, row order -- except row 10 (index 9_CCAACAGA), which
# was filled in last, after all the others.
$skipIndex = 4   # corresponds to row 10 (firstRow + 4)
for ($i = 0; $i -lt $rowCount; $i++) {
    if ($i -ne $skipIndex) {
        $wsP5.Cells.Item($firstRow + $i, 3).Value = $barcodes[$i]
    }
}
$wsP5.Cells.Item($firstRow + $skipIndex, 3).Value = $barcodes[$skipIndex]

# ---------------------------------------------------------------------
# Column B now holds longer text ("...#10_..." .. "...#24_...") -- widen
# it to fit.
# ---------------------------------------------------------------------
$wsP5.Columns.Item(2).ColumnWidth = 35

# ---------------------------------------------------------------------
# Restore each sheet's remembered cursor/selection.
# ---------------------------------------------------------------------
$wsOligo.Range("D9").Select() | Out-Null
$wsTso.Range("E2").Select() | Out-Null
$wsP5.Range("C11").Select() | Out-Null
